$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.217.67"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "2.467.04"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "558.56"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").Value = "163.27"
$ws.Range("E6").Value = "  -1.99%  "
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").Value = "2.467.86"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("E12").Value = "  -4.23%  "
$ws.Range("D13").Value = "4.83"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "2.922.35"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").Value = "69.127.33"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").Value = "23.62"
$ws.Range("D18").Value = "2.432.75"
$ws.Range("E18").Value = "  -3.26%  "
$ws.Range("D19").Value = "10.78"
$ws.Range("E19").Value = "  -4.00%  "
$ws.Range("D20").Value = "343.17"
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("E21").Value = "  -5.48%  "
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").Value = "67.15"
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("D28").Value = "2.596.12"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").Value = "8.20"
$ws.Range("E30").Value = "  -4.98%  "
$ws.Range("E31").Value = "  -6.15%  "
$ws.Range("D32").Value = "7.20"
$ws.Range("E32").Value = "  -5.03%  "
$ws.Range("D33").Value = "439.57"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -4.39%  "
$ws.Range("E36").Value = "  -5.75%  "
$ws.Range("D37").Value = "156.67"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  -3.81%  "
$ws.Range("D41").Value = "17.91"
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("D44").Value = "37.49"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("E45").Value = "  -6.12%  "
$ws.Range("D46").Value = "1.09"
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("E47").Value = "  -4.85%  "
$ws.Range("D48").Value = "133.14"
$ws.Range("E48").Value = "  -4.39%  "
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "0.484"
